$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 47
$ws.Range("F3").Value = 27047
$ws.Range("F5").Value = 648
$ws.Range("F6").Value = 192
$ws.Range("F7").Value = 566
$ws.Range("F9").Value = 376
$ws.Range("F12").Value = 53
$ws.Range("F13").Value = 314
$ws.Range("F14").Value = 96
$ws.Range("F15").Value = 486
$ws.Range("F16").Value = 68
$ws.Range("F17").Value = 1623
$ws.Range("F18").Value = 252
$ws.Range("F19").Value = 697
$ws.Range("F20").Value = 186
$ws.Range("F21").Value = 459
$ws.Range("F23").Value = 108

$ws = $wb.Worksheets.Item(2)
$ws.Range("F2").Value = 4526
$ws.Range("F3").Value = 247
$ws.Range("F8").Value = 42
$ws.Range("F11").Value = 456
$ws.Range("F16").Value = 19
$ws.Range("F17").Value = 75

$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 5190
$ws.Range("F3").Value = 274

$ws = $wb.Worksheets.Item(4)
$ws.Range("F3").Value = 5190
$ws.Range("F4").Value = 274
$ws.Range("F5").Value = 27047
$ws.Range("F6").Value = 4526
$ws.Range("F8").Value = 247
$ws.Range("F9").Value = 648
$ws.Range("F12").Value = 192
$ws.Range("F15").Value = 42
$ws.Range("F18").Value = 456
$ws.Range("F19").Value = 566
$ws.Range("F23").Value = 376
$ws.Range("F26").Value = 53
$ws.Range("F28").Value = 314
$ws.Range("F29").Value = 96
$ws.Range("F31").Value = 19
$ws.Range("F32").Value = 486
$ws.Range("F33").Value = 68
$ws.Range("F34").Value = 75
$ws.Range("F35").Value = 1623
$ws.Range("F36").Value = 253
$ws.Range("F37").Value = 697
$ws.Range("F39").Value = 186
$ws.Range("F40").Value = 459
$ws.Range("F42").Value = 108
